$d = $word.ActiveDocument

# 1) Update the "Updated:" date field's cached result from 2/8/2025 to 3/25/2025
$ok1 = $d.Content.Find.Execute("2/8/2025", $true, $false, $false, $false, $false, $true, 1, $false, "3/25/2025", 2)
if (-not $ok1) { throw "Could not find the date '2/8/2025' to replace." }

# 2) "Displays the active N1/EPR limit computed by the TCI." -> "...N1 or EPR limit..."
$ok2 = $d.Content.Find.Execute("Displays the active N1/EPR limit computed by the TCI.", $true, $false, $false, $false, $false, $true, 1, $false, "Displays the active N1 or EPR limit computed by the TCI.", 2)
if (-not $ok2) { throw "Could not find the N1/EPR limit window description to replace." }

# 3) Reword: "can be used to ensure the engine is operating within safe parameters, which is called"
#    -> "limits are used to ensure the engine is operating within safe parameters. This is called"
$ok3 = $d.Content.Find.Execute("can be used to ensure the engine is operating within safe parameters, which is called", $true, $false, $false, $false, $false, $true, 1, $false, "limits are used to ensure the engine is operating within safe parameters. This is called", 2)
if (-not $ok3) { throw "Could not find the thrust-limit sentence to reword." }
